$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the "22p" capacitor entry. Testing showed C8 and C10 can be
# dropped from that population, so the designator list shrinks from
# "C2, C3, C4, C5, C8, C10" (6 parts) to "C2, C3, C4, C5" (4 parts) and the
# Quantity column (A) drops from 6 to 4 to match.
$ws.Range("B3").Value = "C2, C3, C4, C5"
$ws.Range("A3").Value = 4
$ws.Rows.Item(3).RowHeight = 13.4

# Leave the selection where the editor ended up after making the change.
$ws.Range("B4").Select()
